# Auto-generated Excel COM-interop script applying a scheduled market-data refresh
# to the Tonberry_Profits workbook. Updates currentAveragePrice / Leve profit columns
# (H:N) for the rows whose source values changed; a few rows lose a trailing column
# (the computed "no HQ data" profit cell) when that side of the market disappeared.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 933.44684
$ws.Range("I15").Value = 933.44684
$ws.Range("K15").Value = 2800.34052
$ws.Range("M15").Value = -2631.34052
$ws.Range("H18").Value = 11992.5625
$ws.Range("I18").Value = 6149.625
$ws.Range("J18").Value = 17835.5
$ws.Range("K18").Value = 6149.625
$ws.Range("L18").Value = 17835.5
$ws.Range("M18").Value = -5865.625
$ws.Range("N18").Value = -18403.5
$ws.Range("H28").Value = 2253.5454
$ws.Range("I28").Value = 198.75
$ws.Range("K28").Value = 198.75
$ws.Range("M28").Value = 286.25
$ws.Range("H98").Value = 2752.6
$ws.Range("I98").Value = 2818.0908
$ws.Range("K98").Value = 2818.0908
$ws.Range("M98").Value = -1320.0908
$ws.Range("H122").Value = 2752.6
$ws.Range("I122").Value = 2818.0908
$ws.Range("K122").Value = 8454.2724
$ws.Range("M122").Value = -6004.2724
$ws.Range("H132").Value = 970.3871
$ws.Range("I132").Value = 918.36
$ws.Range("K132").Value = 2755.08
$ws.Range("M132").Value = -225.0799999999999
$ws.Range("H137").Value = 2733.3333
$ws.Range("I137").Value = 1300
$ws.Range("J137").Value = 3450
$ws.Range("K137").Value = 3900
$ws.Range("L137").Value = 10350
$ws.Range("M137").Value = -1350
$ws.Range("N137").Value = -15450
$ws.Range("H141").Value = 919.8
$ws.Range("I141").Value = 919.8
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2759.4
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2420.6
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 370678.94
$ws.Range("I2").Value = 427660.3
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 427660.3
$ws.Range("L2").Value = 300
$ws.Range("M2").Value = -427547.3
$ws.Range("N2").Value = -526
$ws.Range("H32").Value = 3580.0703
$ws.Range("I32").Value = 2876.3582
$ws.Range("J32").Value = 15367.25
$ws.Range("K32").Value = 2876.3582
$ws.Range("L32").Value = 15367.25
$ws.Range("M32").Value = -2589.3582
$ws.Range("N32").Value = -15941.25
$ws.Range("H74").Value = 1701
$ws.Range("I74").Value = 499.2857
$ws.Range("J74").Value = 4505
$ws.Range("K74").Value = 499.2857
$ws.Range("L74").Value = 4505
$ws.Range("M74").Value = 374.7143
$ws.Range("N74").Value = -6253
$ws.Range("H77").Value = 1701
$ws.Range("I77").Value = 499.2857
$ws.Range("J77").Value = 4505
$ws.Range("K77").Value = 2496.4285
$ws.Range("L77").Value = 22525
$ws.Range("M77").Value = 1871.5715
$ws.Range("N77").Value = -31261
$ws.Range("H102").Value = 1508.4286
$ws.Range("I102").Value = 1493.5
$ws.Range("J102").Value = 1514.4
$ws.Range("K102").Value = 1493.5
$ws.Range("L102").Value = 1514.4
$ws.Range("M102").Value = 128.5
$ws.Range("N102").Value = -4758.4
$ws.Range("H110").Value = 357.875
$ws.Range("I110").Value = 373.2857
$ws.Range("K110").Value = 373.2857
$ws.Range("M110").Value = 1671.7143
$ws.Range("H116").Value = 370678.94
$ws.Range("I116").Value = 427660.3
$ws.Range("J116").Value = 300
$ws.Range("K116").Value = 427660.3
$ws.Range("L116").Value = 300
$ws.Range("M116").Value = -425366.3
$ws.Range("N116").Value = -4888

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 370678.94
$ws.Range("I3").Value = 427660.3
$ws.Range("J3").Value = 300
$ws.Range("K3").Value = 427660.3
$ws.Range("L3").Value = 300
$ws.Range("M3").Value = -427546.3
$ws.Range("N3").Value = -528
$ws.Range("H86").Value = 86152.71000000001
$ws.Range("I86").Value = 3551.1333
$ws.Range("K86").Value = 3551.1333
$ws.Range("M86").Value = -2428.1333
$ws.Range("H89").Value = 86152.71000000001
$ws.Range("I89").Value = 3551.1333
$ws.Range("K89").Value = 17755.6665
$ws.Range("M89").Value = -12139.6665
$ws.Range("H105").Value = 2312.4443
$ws.Range("I105").Value = 2057.44
$ws.Range("J105").Value = 5500
$ws.Range("K105").Value = 2057.44
$ws.Range("L105").Value = 5500
$ws.Range("M105").Value = -310.4400000000001
$ws.Range("N105").Value = -8994
$ws.Range("H134").Value = 5292.6294
$ws.Range("I134").Value = 5540.1665
$ws.Range("K134").Value = 16620.4995
$ws.Range("M134").Value = -14085.4995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2045.6666
$ws.Range("I31").Value = 1832.8334
$ws.Range("J31").Value = 2471.3333
$ws.Range("K31").Value = 1832.8334
$ws.Range("L31").Value = 2471.3333
$ws.Range("M31").Value = -1537.8334
$ws.Range("N31").Value = -3061.3333
$ws.Range("H34").Value = 2045.6666
$ws.Range("I34").Value = 1832.8334
$ws.Range("J34").Value = 2471.3333
$ws.Range("K34").Value = 1832.8334
$ws.Range("L34").Value = 2471.3333
$ws.Range("M34").Value = -1630.8334
$ws.Range("N34").Value = -2875.3333
$ws.Range("H58").Value = 2072380.6
$ws.Range("I58").Value = 2719225.5
$ws.Range("J58").Value = 2477
$ws.Range("K58").Value = 2719225.5
$ws.Range("L58").Value = 2477
$ws.Range("M58").Value = -2719022.5
$ws.Range("N58").Value = -2883
$ws.Range("H105").Value = 892.5833
$ws.Range("I105").Value = 879.5
$ws.Range("K105").Value = 879.5
$ws.Range("M105").Value = 867.5
$ws.Range("H132").Value = 2360.0833
$ws.Range("I132").Value = 1290.375
$ws.Range("K132").Value = 3871.125
$ws.Range("M132").Value = -1341.125
$ws.Range("H136").Value = 2072380.6
$ws.Range("I136").Value = 2719225.5
$ws.Range("J136").Value = 2477
$ws.Range("K136").Value = 8157676.5
$ws.Range("L136").Value = 7431
$ws.Range("M136").Value = -8155126.5
$ws.Range("N136").Value = -12531

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 599.6667
$ws.Range("I11").Value = 519.8
$ws.Range("K11").Value = 1559.4
$ws.Range("M11").Value = -1419.4
$ws.Range("H25").Value = 1949.8
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 1949.8
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 5849.4
$ws.Range("N25").Value = -6187.4
$ws.Range("M25").ClearContents()
$ws.Range("H30").Value = 1949.8
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 1949.8
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 5849.4
$ws.Range("N30").Value = -6053.4
$ws.Range("M30").ClearContents()
$ws.Range("H32").Value = 2416.3333
$ws.Range("J32").Value = 2416.3333
$ws.Range("L32").Value = 7248.999899999999
$ws.Range("N32").Value = -7814.999899999999
$ws.Range("H34").Value = 9484.546
$ws.Range("I34").Value = 11236.667
$ws.Range("J34").Value = 1600
$ws.Range("K34").Value = 33710.001
$ws.Range("L34").Value = 4800
$ws.Range("M34").Value = -33626.001
$ws.Range("N34").Value = -4968
$ws.Range("H131").Value = 15660.745
$ws.Range("J131").Value = 17445.834
$ws.Range("L131").Value = 52337.50199999999
$ws.Range("N131").Value = -62417.50199999999
$ws.Range("H140").Value = 2412.3225
$ws.Range("I140").Value = 952.2353000000001
$ws.Range("J140").Value = 4185.2856
$ws.Range("K140").Value = 2856.7059
$ws.Range("L140").Value = 12555.8568
$ws.Range("M140").Value = 2323.2941
$ws.Range("N140").Value = -22915.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3320
$ws.Range("I80").Value = 3490
$ws.Range("K80").Value = 3490
$ws.Range("M80").Value = -2492
$ws.Range("H83").Value = 3320
$ws.Range("I83").Value = 3490
$ws.Range("K83").Value = 17450
$ws.Range("M83").Value = -12458
$ws.Range("H132").Value = 1750598.2
$ws.Range("I132").Value = 2026166.5
$ws.Range("K132").Value = 6078499.5
$ws.Range("M132").Value = -6075969.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H100").Value = 1240.1666
$ws.Range("I100").Value = 1088.4
$ws.Range("K100").Value = 2176.8
$ws.Range("M100").Value = -1635.8
$ws.Range("H132").Value = 2462.5881
$ws.Range("I132").Value = 1686.7693
$ws.Range("K132").Value = 5060.3079
$ws.Range("M132").Value = -2530.3079
